# "improved canala ponding: distributed water throughout cell"
#
# Adds two new rows to the file-pointers table describing the new
# pointers used by the pickled initial-zeta raster and the template
# output raster, widens column A slightly so the new, longer label
# ("template_output_raster") fits, and leaves the selection where the
# author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11: initial_zeta_pickle -> data/new_area/best_initial_zeta.p
$ws.Range("A11").Value = "initial_zeta_pickle"
$ws.Range("A12").Value = "template_output_raster"
$ws.Range("B11").Value = "data/new_area/best_initial_zeta.p"
# New row 12: template_output_raster -> data/dtm_depth_padded.tif (reuses the existing DEM path string)
$ws.Range("B12").Value = "data/dtm_depth_padded.tif"

# Widen column A (20 -> ~24.14 characters) so the longer labels fit.
$ws.Columns.Item(1).ColumnWidth = 23.333333333333332

# Match the author's final selection/scroll position.
$selectResult = $ws.Range("N19").Select()

# Best-effort: restore the saved window geometry too.
$win = $excel.ActiveWindow
$win.Left = 5145
$win.Top = 5610
$win.Width = 28800
$win.Height = 15435
